$d = $word.ActiveDocument

# Remove the now-unused "Abstract Title" custom paragraph style entirely.
$abstractTitle = $d.Styles("AbstractTitle")
$abstractTitle.Delete()

# The "Abstract" style's paragraph spacing-before changes from 100 (5pt)
# to 300 (15pt) twips, matching spacing-after.
$abstract = $d.Styles("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 15
